# repull data, push all data, mean calculation
# Update the dSF (column F) values for the affected rows to reflect the
# re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 6
    4  = -3
    5  = 1
    6  = -2
    9  = -3
    10 = 1
    12 = -3
    14 = 1
    15 = 1
    18 = -5
    19 = 1
    20 = 8
    21 = -2
    22 = 2
    23 = -2
    24 = -4
    25 = -2
    26 = 2
    27 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
